$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving new text values. NumberFormat is forced to Text ("@")
# before assignment so Excel does not auto-coerce numeric-looking strings
# (e.g. "1.006", "14.80", "8.440") into floating point numbers, which would
# silently drop significant trailing zeros / change the stored representation.
$targetCells = @(
    'D2',
    'E2',
    'E3',
    'E4',
    'D5',
    'E5',
    'D6',
    'E6',
    'D7',
    'E7',
    'D8',
    'E8',
    'D9',
    'E9',
    'D10',
    'E10',
    'D11',
    'E11',
    'D12',
    'E12',
    'D13',
    'E13',
    'D14',
    'E14',
    'E15',
    'D16',
    'E16',
    'D17',
    'E17',
    'D18',
    'E18',
    'D19',
    'E19',
    'D20',
    'E20',
    'D21',
    'E21',
    'D22',
    'E22',
    'D23',
    'E23',
    'E24',
    'D25',
    'E25',
    'D26',
    'E26',
    'E27',
    'D28',
    'E28',
    'D29',
    'E29',
    'E30',
    'D31',
    'E31',
    'D32',
    'E32',
    'E33',
    'B34',
    'C34',
    'D34',
    'E34',
    'B35',
    'C35',
    'D35',
    'E35',
    'D36',
    'E36',
    'D37',
    'E37',
    'D38',
    'E38',
    'D39',
    'E39',
    'D40',
    'E40',
    'D41',
    'E41',
    'D42',
    'E42',
    'E43',
    'D44',
    'E44',
    'E45',
    'D46',
    'E46',
    'D47',
    'E47',
    'D48',
    'E48',
    'D49',
    'E49',
    'D50',
    'E50',
    'D51',
    'E51'
)

foreach ($ref in $targetCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '28.053.58'
$ws.Range('E2').Value = '  +1.91%  '
$ws.Range('E3').Value = '  +1.76%  '
$ws.Range('E4').Value = '  -0.86%  '
$ws.Range('D5').Value = '315.81'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('D7').Value = '0.4825'
$ws.Range('E7').Value = '  +0.81%  '
$ws.Range('D8').Value = '0.3817'
$ws.Range('E8').Value = '  +1.06%  '
$ws.Range('D9').Value = '0.07359'
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').Value = '0.9323'
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('D11').Value = '20.74'
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('D12').Value = '0.07822'
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('D13').Value = '1.884.57'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').Value = '5.508'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').Value = '91.16'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('D18').Value = '0.000008816'
$ws.Range('E18').Value = '  -1.46%  '
$ws.Range('D19').Value = '1.006'
$ws.Range('E19').Value = '  -0.77%  '
$ws.Range('D20').Value = '28.059.18'
$ws.Range('E20').Value = '  +1.72%  '
$ws.Range('D21').Value = '14.80'
$ws.Range('E21').Value = '  -0.82%  '
$ws.Range('D22').Value = '5.152'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').Value = '2.137.74'
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('D25').Value = '156.51'
$ws.Range('E25').Value = '  +1.90%  '
$ws.Range('D26').Value = '1.923'
$ws.Range('E26').Value = '  -1.92%  '
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').Value = '2.101'
$ws.Range('E28').Value = '  +3.96%  '
$ws.Range('D29').Value = '116.34'
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('E30').Value = '  -1.15%  '
$ws.Range('D31').Value = '0.08904'
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('D32').Value = '3.355'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('E33').Value = '  +2.33%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '4.672'
$ws.Range('E34').Value = '  +1.31%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '0.7647'
$ws.Range('E35').Value = '  +1.90%  '
$ws.Range('D36').Value = '2.613'
$ws.Range('E36').Value = '  -2.95%  '
$ws.Range('D37').Value = '0.02043'
$ws.Range('E37').Value = '  -1.52%  '
$ws.Range('D38').Value = '1.097'
$ws.Range('E38').Value = '  -1.92%  '
$ws.Range('D39').Value = '0.05291'
$ws.Range('E39').Value = '  -0.28%  '
$ws.Range('D40').Value = '0.5476'
$ws.Range('E40').Value = '  +2.28%  '
$ws.Range('D41').Value = '2.981'
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('D42').Value = '7.003'
$ws.Range('E42').Value = '  -1.48%  '
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('D44').Value = '8.440'
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('D46').Value = '0.4825'
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('D47').Value = '107.17'
$ws.Range('E47').Value = '  +3.85%  '
$ws.Range('D48').Value = '1.006'
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('D49').Value = '1.652'
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('D50').Value = '68.13'
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('D51').Value = '0.06098'
$ws.Range('E51').Value = '  +0.00%  '

# Restore the default cell style so the cells keep no explicit style index,
# matching the original workbook formatting (only the text content changed).
foreach ($ref in $targetCells) {
    $ws.Range($ref).Style = "Normal"
}
